$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("carrito")

# Update the search term (B2): "reloj hombre" -> "Lapicero"
$ws.Range("B2").Value = "Lapicero"

# Clear the now-unused "producto"/"precio" columns (C1:D2) - content only,
# D2 keeps its numeric cell/style but becomes empty.
$ws.Range("C1:D2").ClearContents()

# Column B widens (new text no longer fits the old width).
$ws.Columns("B").ColumnWidth = 15.8

# Move the active selection to B3.
$ws.Range("B3").Select() | Out-Null
